$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 2 (the first data row after the header). This shifts all
# subsequent rows (3..69) up by one, matching the target diff where
# the old row N becomes row N-1 for N = 3..69.
$ws.Rows.Item(2).Delete()
